$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "fullName"
$ws.Range("B1").Value = "email"
$ws.Range("C1").Value = "currentAddress"
$ws.Range("D1").Value = "permanentAdd"

# Data rows
$ws.Range("A2").Value = "Davinder Sharma"
$ws.Range("B2").Value = "dav@gmail.com"
$ws.Range("C2").Value = "Bengaluru"
$ws.Range("D2").Value = "patna"

$ws.Range("A3").Value = "Rohit"
$ws.Range("B3").Value = "rohit@gmail.cm"
$ws.Range("C3").Value = "marathali"
$ws.Range("D3").Value = "belgam"

$ws.Range("A4").Value = "Thushar"
$ws.Range("B4").Value = "thushar@gmail.com"
$ws.Range("C4").Value = "multiplex"
$ws.Range("D4").Value = "coorg"

$ws.Range("A5").Value = "Soumik"
$ws.Range("B5").Value = "soumik@gmail.com"
$ws.Range("C5").Value = "vkrPG"
$ws.Range("D5").Value = "aasam"

# Hyperlinks on email column
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:dav@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:rohit@gmail.cm")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:thushar@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:soumik@gmail.com")
